# Weekly refresh of the price data: rows 2-12 get their D/J/K/L/M/O/P
# values shuffled among themselves (row 9 keeps its own values).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row N gets the D/J/K/L/M/O/P values that used to live on old row Source(N)
$rowSource = @{
    2  = 5
    3  = 6
    4  = 12
    5  = 8
    6  = 4
    7  = 11
    8  = 10
    9  = 9
    10 = 7
    11 = 2
    12 = 3
}

# Snapshot the "before" values for the columns that move, keyed by row.
$cols = @("D", "J", "K", "L", "M", "O", "P")
$snapshot = @{}
foreach ($r in 2..12) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# Write the new values, pulled from the snapshot of the source row.
foreach ($r in 2..12) {
    $src = $rowSource[$r]
    $srcVals = $snapshot[$src]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = $srcVals[$c]
    }
}
